$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the two newly tracked "linear regression" row entries (row 15 = DS intro, row 16 = linear regression)
$ws.Range("D15").Value = 0.5
$ws.Range("D16").Value = 1.5

# Update the active selection to reflect where the user left off (D17)
$ws.Range("D17").Select()
